# Weekly update for "Fruta, Terminal Hortofrutícola Agro Chillán - Cereza":
# a new week's record is inserted at the top of the data block (row 5),
# pushing every existing record down by one row, and the final (oldest)
# record that falls off the bottom is re-appended as the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5; this shifts the existing rows 5..59
# down to 6..60 and extends the used range to A1:T60, matching the rest
# of the constant (non-varying) columns for every record in this sheet.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with this week's record. The
# "constant" columns (A,B,C,E,F,G,H,I,J,T) are identical for every row in
# this sheet, and K/L (Variedad/Calidad) stay "Lapins"/"Primera" here.
$ws.Cells.Item(5, 1).Value  = 7
$ws.Cells.Item(5, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5, 3).Value  = "Ñuble"
$ws.Cells.Item(5, 4).Value  = 44530
$ws.Cells.Item(5, 5).Value  = 16
$ws.Cells.Item(5, 6).Value  = "Fruta"
$ws.Cells.Item(5, 7).Value  = 100103
$ws.Cells.Item(5, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(5, 9).Value  = 100103001
$ws.Cells.Item(5, 10).Value = "Cereza"
$ws.Cells.Item(5, 11).Value = "Lapins"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 60
$ws.Cells.Item(5, 14).Value = 15000
$ws.Cells.Item(5, 15).Value = 16000
$ws.Cells.Item(5, 16).Value = 15500
$ws.Cells.Item(5, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(5, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(5, 19).Value = 1550
$ws.Cells.Item(5, 20).Value = 10
